$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.731.73"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "3.436.30"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'573.23"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("D6").Value = "'158.87"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.438.20"
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("D9").Value = "'0.584"
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("E11").Value = "  -2.99%  "
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").Value = "4.033.69"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "'27.60"
$ws.Range("E15").Value = "  -3.96%  "
$ws.Range("D16").Value = "'0.0000179"
$ws.Range("E16").Value = "  -7.57%  "
$ws.Range("D17").Value = "64.790.81"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").Value = "3.446.42"
$ws.Range("E18").Value = "  -1.78%  "
$ws.Range("D19").Value = "'6.29"
$ws.Range("E19").Value = "  -2.92%  "
$ws.Range("D20").Value = "'13.87"
$ws.Range("E20").Value = "  -3.72%  "
$ws.Range("D21").Value = "'380.23"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").Value = "'7.94"
$ws.Range("E22").Value = "  -4.60%  "
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").Value = "'72.05"
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("E26").Value = "  -4.11%  "
$ws.Range("D27").Value = "'9.83"
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("E31").Value = "  -2.67%  "
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("D33").Value = "'23.19"
$ws.Range("E33").Value = "  -2.03%  "
$ws.Range("D34").Value = "'7.00"
$ws.Range("E34").Value = "  -4.54%  "
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("D36").Value = "'161.37"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("D38").Value = "2.891.66"
$ws.Range("E38").Value = "  -4.18%  "
$ws.Range("D39").Value = "'0.0743"
$ws.Range("E39").Value = "  -4.63%  "
$ws.Range("D40").Value = "'26.21"
$ws.Range("E40").Value = "  -3.92%  "
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").Value = "'4.54"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("D43").Value = "'42.90"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  -3.05%  "
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").Value = "'25.97"
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("D47").Value = "'2.26"
$ws.Range("E47").Value = "  +1.92%  "
$ws.Range("E48").Value = "  -3.28%  "
$ws.Range("D49").Value = "'316.12"
$ws.Range("E49").Value = "  -3.11%  "
$ws.Range("D50").Value = "'6.50"
$ws.Range("E50").Value = "  -3.84%  "
$ws.Range("E51").Value = "  -3.33%  "
